# Updates the coin symbol/price/volume table to the latest scraped snapshot.
# (commit: "Updated symbol list on Fri Feb  3 14:34:24 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link columns (plain text, rows 18-23 shifted by one new listing) ---
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B20").Value = "LEO"
$ws.Range("C20").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"

# --- Price / Volume(1h) columns ---
# These columns store plain numeric-looking text (e.g. "330.58", "0.69%").
# Force the cell format to Text *before* assigning so Excel keeps the exact
# string (incl. trailing zeros / percent signs) instead of coercing to a number.
$priceVolumeCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D10",
    "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16",
    "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23",
    "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D39", "E39", "D40", "E40", "D41", "D42",
    "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48",
    "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($cellRef in $priceVolumeCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "330.58"
$ws.Range("E2").Value = "0.69%"
$ws.Range("D3").Value = "39.11"
$ws.Range("E3").Value = "-2.85%"
$ws.Range("D4").Value = "5.667"
$ws.Range("E4").Value = "1.01%"
$ws.Range("D5").Value = "0.08023"
$ws.Range("E5").Value = "-1.57%"
$ws.Range("D6").Value = "4.485"
$ws.Range("E6").Value = "-1.45%"
$ws.Range("D7").Value = "8.600"
$ws.Range("E7").Value = "-0.82%"
$ws.Range("D8").Value = "1.941"
$ws.Range("E8").Value = "-1.68%"
$ws.Range("D10").Value = "0.9191"
$ws.Range("E10").Value = "-3.14%"
$ws.Range("D11").Value = "0.1238"
$ws.Range("E11").Value = "-3.14%"
$ws.Range("D12").Value = "0.1938"
$ws.Range("E12").Value = "-2.31%"
$ws.Range("D13").Value = "8.698"
$ws.Range("E13").Value = "16.28%"
$ws.Range("D14").Value = "0.09194"
$ws.Range("E14").Value = "-0.35%"
$ws.Range("D15").Value = "0.03508"
$ws.Range("E15").Value = "-1.53%"
$ws.Range("D16").Value = "0.1049"
$ws.Range("E16").Value = "9.03%"
$ws.Range("D17").Value = "0.001315"
$ws.Range("E17").Value = "-0.32%"
$ws.Range("D18").Value = "0.04440"
$ws.Range("E18").Value = "0.31%"
$ws.Range("D19").Value = "0.006329"
$ws.Range("E19").Value = "-1.04%"
$ws.Range("D20").Value = "3.363"
$ws.Range("E20").Value = "-0.28%"
$ws.Range("D21").Value = "0.3457"
$ws.Range("E21").Value = "-1.77%"
$ws.Range("D22").Value = "0.1373"
$ws.Range("E22").Value = "-1.67%"
$ws.Range("D23").Value = "0.2697"
$ws.Range("E23").Value = "8.71%"
$ws.Range("D24").Value = "0.001255"
$ws.Range("E24").Value = "-0.12%"
$ws.Range("D25").Value = "0.004497"
$ws.Range("E25").Value = "4.67%"
$ws.Range("D26").Value = "0.0001202"
$ws.Range("E26").Value = "1.18%"
$ws.Range("D39").Value = "0.02551"
$ws.Range("E39").Value = "0.54%"
$ws.Range("D40").Value = "0.05437"
$ws.Range("E40").Value = "4.23%"
$ws.Range("D41").Value = "0.007538"
$ws.Range("D42").Value = "0.009915"
$ws.Range("E42").Value = "10.26%"
$ws.Range("D43").Value = "0.1401"
$ws.Range("E43").Value = "-2.53%"
$ws.Range("D44").Value = "0.002111"
$ws.Range("E44").Value = "-3.46%"
$ws.Range("D45").Value = "0.01154"
$ws.Range("E45").Value = "16.40%"
$ws.Range("D46").Value = "0.00006803"
$ws.Range("E46").Value = "1.74%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.33%"
$ws.Range("D48").Value = "0.003056"
$ws.Range("E48").Value = "6.53%"
$ws.Range("D49").Value = "0.002283"
$ws.Range("E49").Value = "-0.72%"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.33%"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "0.33%"
